$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 724 (the "好奇心旺盛なホッキョクギツネ" post) entirely, shifting all
# subsequent rows up by one.
$ws.Rows(724).Delete()
